$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 793.125
$ws.Range("J29").Value = 1041.6666
$ws.Range("L29").Value = 3124.9998
$ws.Range("N29").Value = -3686.9998
$ws.Range("H38").Value = 163.61539
$ws.Range("I38").Value = 163.61539
$ws.Range("K38").Value = 490.84617
$ws.Range("M38").Value = -118.84617
$ws.Range("H40").Value = 2062.5
$ws.Range("J40").Value = 6000
$ws.Range("L40").Value = 6000
$ws.Range("N40").Value = -6350
$ws.Range("H87").Value = 49997.5
$ws.Range("J87").Value = 49997.5
$ws.Range("L87").Value = 49997.5
$ws.Range("N87").Value = -52493.5
$ws.Range("H90").Value = 49997.5
$ws.Range("J90").Value = 49997.5
$ws.Range("L90").Value = 149992.5
$ws.Range("N90").Value = -162472.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2288.5
$ws.Range("I45").Value = 2577
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2577
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -2200
$ws.Range("N45").Value = -2754
$ws.Range("H122").Value = 2308
$ws.Range("I122").Value = 2309.6
$ws.Range("K122").Value = 6928.799999999999
$ws.Range("M122").Value = -4478.799999999999
$ws.Range("H132").Value = 3384.72
$ws.Range("I132").Value = 3148.4285
$ws.Range("J132").Value = 4625.25
$ws.Range("K132").Value = 9445.2855
$ws.Range("L132").Value = 13875.75
$ws.Range("M132").Value = -6915.2855
$ws.Range("N132").Value = -18935.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 182
$ws.Range("I22").Value = 182
$ws.Range("K22").Value = 182
$ws.Range("M22").Value = -9
$ws.Range("H81").Value = 45159.832
$ws.Range("J81").Value = 45159.832
$ws.Range("L81").Value = 45159.832
$ws.Range("N81").Value = -47281.832
$ws.Range("H84").Value = 45159.832
$ws.Range("J84").Value = 45159.832
$ws.Range("L84").Value = 135479.496
$ws.Range("N84").Value = -146087.496
$ws.Range("H99").Value = 2977.8333
$ws.Range("I99").Value = 2749.111
$ws.Range("J99").Value = 3664
$ws.Range("K99").Value = 2749.111
$ws.Range("L99").Value = 3664
$ws.Range("M99").Value = -1251.111
$ws.Range("N99").Value = -6660
$ws.Range("H134").Value = 9044.6
$ws.Range("I134").Value = 9534.857
$ws.Range("K134").Value = 28604.571
$ws.Range("M134").Value = -26069.571
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1385.2222
$ws.Range("I31").Value = 1224.5
$ws.Range("J31").Value = 1706.6666
$ws.Range("K31").Value = 1224.5
$ws.Range("L31").Value = 1706.6666
$ws.Range("M31").Value = -929.5
$ws.Range("N31").Value = -2296.6666
$ws.Range("H34").Value = 1385.2222
$ws.Range("I34").Value = 1224.5
$ws.Range("J34").Value = 1706.6666
$ws.Range("K34").Value = 1224.5
$ws.Range("L34").Value = 1706.6666
$ws.Range("M34").Value = -1022.5
$ws.Range("N34").Value = -2110.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1516.2
$ws.Range("I3").Value = 1516.2
$ws.Range("K3").Value = 4548.6
$ws.Range("M3").Value = -4436.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2126998.2
$ws.Range("J70").Value = 7994
$ws.Range("L70").Value = 7994
$ws.Range("N70").Value = -8534
$ws.Range("H73").Value = 2126998.2
$ws.Range("J73").Value = 7994
$ws.Range("L73").Value = 7994
$ws.Range("N73").Value = -9866
$ws.Range("H132").Value = 1823.4445
$ws.Range("I132").Value = 1689.375
$ws.Range("K132").Value = 5068.125
$ws.Range("M132").Value = -2538.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 21087.846
$ws.Range("I7").Value = 21087.846
$ws.Range("K7").Value = 21087.846
$ws.Range("M7").Value = -20975.846
$ws.Range("H40").Value = 7980.4
$ws.Range("I40").Value = 6633.6665
$ws.Range("K40").Value = 6633.6665
$ws.Range("M40").Value = -6497.6665
$ws.Range("H46").Value = 5463.6665
$ws.Range("I46").Value = 5608.5625
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 5608.5625
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -5420.5625
$ws.Range("N46").Value = -5376
$ws.Range("H55").Value = 341.33334
$ws.Range("I55").Value = 323.75
$ws.Range("J55").Value = 355.4
$ws.Range("K55").Value = 323.75
$ws.Range("L55").Value = 355.4
$ws.Range("M55").Value = -150.75
$ws.Range("N55").Value = -701.4
$ws.Range("H126").Value = 21087.846
$ws.Range("I126").Value = 21087.846
$ws.Range("K126").Value = 63263.538
$ws.Range("M126").Value = -60793.538
$ws.Range("H127").Value = 44000
$ws.Range("J127").Value = 44000
$ws.Range("L127").Value = 44000
$ws.Range("N127").Value = -53920
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("H136").Value = 2912.4614
$ws.Range("I136").Value = 2386.2
$ws.Range("K136").Value = 7158.599999999999
$ws.Range("M136").Value = -4608.599999999999
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6499.75
$ws.Range("I81").Value = 4000
$ws.Range("J81").Value = 8999.5
$ws.Range("K81").Value = 8000
$ws.Range("L81").Value = 17999
$ws.Range("M81").Value = -6939
$ws.Range("N81").Value = -20121
$ws.Range("H84").Value = 6499.75
$ws.Range("I84").Value = 4000
$ws.Range("J84").Value = 8999.5
$ws.Range("K84").Value = 40000
$ws.Range("L84").Value = 89995
$ws.Range("M84").Value = -34696
$ws.Range("N84").Value = -100603
$ws.Range("H122").Value = 2899.125
$ws.Range("I122").Value = 2798.5
$ws.Range("K122").Value = 8395.5
$ws.Range("M122").Value = -5945.5
$ws.Range("H132").Value = 2945.2
$ws.Range("I132").Value = 2575.3333
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 7725.999899999999
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -5195.999899999999
$ws.Range("N132").Value = -15560
